# ---------------------------------------------------------------------------
# edit.ps1 - applies the "Added Compilation notes to the Report" change set
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# ===========================================================================
# Part 1: "...com o numero identificador do erro (" ->
#         "...com o numero negativo identificador do erro ("
# ===========================================================================
$find = $d.Content.Find
$find.ClearFormatting()
$ok = $find.Execute(
    "número identificador do erro (",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "número negativo identificador do erro (",
    2)
Write-Output "Part1 replace ok: $ok"

# ===========================================================================
# Part 2a: append a new sentence after "...a mensagem "SERVER CLOSED"."
#          (the trailing run of that paragraph is rebuilt, together with a
#          brand new run, via InsertXML over a range that spans through the
#          end of the paragraph, so the splice fully replaces that tail
#          instead of merely inserting next to it)
# ===========================================================================
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute(
    " a mensagem " + [char]8220 + "SERVER CLOSED" + [char]8221 + ".",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
Write-Output "Part2a find ok: $found2"
if ($found2) {
    $tailRun = $find2.Parent
    $ownerPara = $tailRun.Paragraphs.Item(1)
    $spliceRange = $d.Range($tailRun.Start, $ownerPara.Range.End)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="228" w:after="228"/><w:jc w:val="both"/>' + `
        '<w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>' + `
        '<w:t xml:space="preserve"> a mensagem &#x201c;SERVER CLOSED&#x201d;. </w:t></w:r>' + `
        '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>' + `
        '<w:t>Todos os descritores de ficheiros abertos são fechados</w:t></w:r>' + `
        '</w:p>'
    $spliceRange.InsertXML($xml)
}

# ===========================================================================
# Part 2b: merge the bookmark/page-break paragraph with the heading
#          paragraph that follows it, dropping the page break and fixing
#          up the bookmarks.
# ===========================================================================
$pBookmarkIdx = -1
$pHeadingIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "Mecanismos de Sincroniza") {
        $pHeadingIdx = $i
        $pBookmarkIdx = $i - 1
        break
    }
}
Write-Output "bookmarkIdx=$pBookmarkIdx headingIdx=$pHeadingIdx"

if ($pBookmarkIdx -gt 0) {
    $pB = $d.Paragraphs.Item($pBookmarkIdx)
    $pH = $d.Paragraphs.Item($pHeadingIdx)
    $mergeRange = $d.Range($pB.Range.Start, $pH.Range.End)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr/></w:pPr><w:bookmarkStart w:id="1" w:name="__DdeLink__0_1784306064"/><w:bookmarkStart w:id="2" w:name="__DdeLink__9_512900471"/><w:bookmarkEnd w:id="1"/><w:bookmarkEnd w:id="2"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Mecanismos de Sincronização Utilizados</w:t></w:r></w:p>'
    $mergeRange.InsertXML($xml)
}

# ===========================================================================
# Part 3: append the new "Notas relativas a Compilacao" section at the end
#         of the document body.
# ===========================================================================
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$emptyPara = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:r></w:p>'

$headingPara = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Notas relativas à Compilação</w:t></w:r></w:p>'

$bodyPara = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="171" w:after="171"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">Cada programa tem um </w:t></w:r>' + `
    '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:iCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>makefile</w:t></w:r>' + `
    '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> próprio. Para compilar a totalidade do programa (client + server), deve-se correr o script “build.sh”, que executa o </w:t></w:r>' + `
    '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i/><w:iCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">makefile </w:t></w:r>' + `
    '<w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>de cada programa e coloca o binário resultante de cada programa na mesma pasta (condição necessária para o funcionamento).</w:t></w:r>' + `
    '</w:p>'

$fullXml = $emptyPara + $emptyPara + $emptyPara + $headingPara + $bodyPara

$endRange.InsertXML($fullXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
